$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value
$ws.Range("A2").Value = "rrrrrrrrr"

# Remove row 3 entirely (shrinks used range to A1:F2)
$ws.Rows(3).Delete()
